$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plant Data")

# --- Period 1 block (columns A-C): add Carbon Intensity (CI) column C ---
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0.5
$ws.Range("C4").Value = 0.8
$ws.Range("C5").Value = 1
$ws.Range("C2:C5").HorizontalAlignment = -4108

# --- Period 2 block (columns D-F): add Carbon Intensity (CI) column F ---
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0.5
$ws.Range("F4").Value = 0.8
$ws.Range("F5").Value = 1
$ws.Range("F2:F5").HorizontalAlignment = -4108

# --- Period 3 block (columns G-I): add Carbon Intensity (CI) column I ---
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0.5
$ws.Range("I4").Value = 0.8
$ws.Range("I5").Value = 1
$ws.Range("I2:I5").HorizontalAlignment = -4108

# --- NET_CI totals on row 6, matching the formatting already used by the
#     adjacent TOTAL cells (copy format first, then set formula/number format) ---
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Formula = "=(B2*C2)+(B3*C3)+(B4*C4)+(B5*C5)"

$ws.Range("B6").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Formula = "=(E2*F2)+(E3*F3)+(E4*F4)+(E5*F5)"
$ws.Range("F6").NumberFormat = "0.00"

$ws.Range("B6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").Formula = "=(H2*I2)+(H3*I3)+(H4*I4)+(H5*I5)"
$ws.Range("I6").NumberFormat = "0.00"

$excel.CutCopyMode = $false

# --- Cosmetic: move the live selection like the saved workbook did ---
$ws.Range("E12").Select()
